$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell reference -> new text value for the "Updated symbol list" data refresh.
$updates = [ordered]@{
    'E2' = '1.16%'
    'D3' = '27.12'
    'E3' = '0.96%'
    'E4' = '1.32%'
    'D5' = '0.06184'
    'E5' = '3.46%'
    'D6' = '6.689'
    'E6' = '0.62%'
    'D7' = '0.8506'
    'E7' = '-0.69%'
    'D8' = '0.9138'
    'E8' = '-0.75%'
    'D9' = '0.1406'
    'E9' = '1.28%'
    'D10' = '0.04664'
    'E10' = '-5.07%'
    'E11' = '0.94%'
    'D12' = '0.03155'
    'E12' = '3.50%'
    'D13' = '0.09044'
    'E13' = '-0.89%'
    'D14' = '0.001543'
    'E14' = '0.09%'
    'D15' = '0.0006171'
    'E15' = '2.21%'
    'D16' = '0.006000'
    'E16' = '-2.77%'
    'D17' = '3.458'
    'D18' = '3.181'
    'E18' = '1.09%'
    'D19' = '2.176'
    'E19' = '0.36%'
    'E20' = '-0.96%'
    'E21' = '0.86%'
    'D22' = '4.092'
    'E22' = '1.34%'
    'D23' = '0.04236'
    'E23' = '0.08%'
    'D24' = '0.001215'
    'E24' = '-0.18%'
    'E25' = '-5.53%'
    'E26' = '0.04%'
    'D40' = '0.03898'
    'E40' = '1.85%'
    'D41' = '0.1111'
    'E41' = '-0.07%'
    'D42' = '0.004098'
    'E42' = '8.85%'
    'B43' = 'CEJI'
    'C43' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'D43' = '0.002184'
    'E43' = '-10.12%'
    'B44' = 'LocalTraders'
    'C44' = 'https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct'
    'D44' = '0.01390'
    'E44' = '-8.00%'
    'D45' = '0.00005158'
    'E45' = '1.09%'
    'E46' = '0.04%'
    'E48' = '28.57%'
    'E49' = '0.04%'
    'E50' = '0.04%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Force text storage so numeric- and percent-looking strings (e.g. "1.16%", "27.12")
    # are kept as literal text, matching the workbook's existing inline-string cells
    # instead of being auto-converted into numbers/percentages by Excel.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
